$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.918.99"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.903.57"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'568.56"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").Value = "'143.56"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.500"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.901.41"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'6.99"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "'0.429"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'32.57"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "3.387.57"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "61.918.19"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "2.910.28"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "'429.33"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'13.01"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").Value = "'6.86"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'78.78"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'11.96"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "'10.23"
$ws.Range("E26").Value = "  -7.52%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'2.02"
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("E29").Value = "  +12.04%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").Value = "'0.953"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").Value = "'5.38"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").Value = "'48.80"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").Value = "'2.88"
$ws.Range("E39").Value = "  -6.84%  "
$ws.Range("D40").Value = "'1.90"
$ws.Range("E40").Value = "  -4.87%  "
$ws.Range("D42").Value = "'41.16"
$ws.Range("E42").Value = "  +5.77%  "
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "2.708.49"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "'133.35"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "'0.0337"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'349.01"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "'0.000208"
$ws.Range("E51").Value = "  +10.29%  "
